$d = $word.ActiveDocument

# --- Fix the "logginf" typo -> "login" ------------------------------------
# Locate the mistyped word and the offending character ("f" at the 7th
# position of "logginf").
$text = $d.Content.Text
$idx = $text.IndexOf("logginf")
if ($idx -ge 0) {
    $fPos = $idx + 6

    # Replace the "f" with "g", turning "logginf" into "login".
    $charRange = $d.Range($fPos, $fPos + 1)
    $charRange.Text = "g"

    # Force a run-split between "loggin" and the corrected "g" by
    # temporarily bookmarking the freshly typed character...
    $wrapRange = $d.Range($fPos, $fPos + 1)
    $d.Bookmarks.Add("_GoBack", $wrapRange)

    # ...then collapse/move the bookmark to sit right after the "g",
    # marking the end of the last edit, as Word does automatically.
    $collapsedRange = $d.Range($fPos + 1, $fPos + 1)
    $d.Bookmarks.Add("_GoBack", $collapsedRange)
}
